# EI Variable Installments T2 scenarios
# - Insert a new "waittopageload1" step (row) into the
#   "Edit Repayment Schedule" sheet's step list, keeping the green
#   amount-cell formatting used by the similar "waittopageload" row.
# - Make "Edit Repayment Schedule" the active/selected sheet (it was
#   "NewLoanInput" before), with the in-sheet selection sitting on the
#   first cell of the newly shifted "clickonsubmit" step.

$wb = $excel.ActiveWorkbook

$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row above the existing row 6 ("clickonsubmit"/"Submit"),
# which pushes every following step down by one.
$wsEdit.Rows.Item(6).Insert()

$wsEdit.Range("A6").Value = "waittopageload1"
$wsEdit.Range("B6").Value = 2000

# The inserted row inherits row 5's (text) formatting by default; match
# the numeric "amount" styling used by the other wait/adjust steps
# (e.g. B3/B4) by copying B3's format onto the new B6.
$wsEdit.Range("B3").Copy()
[void]$wsEdit.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the in-sheet selection onto the (now shifted) "clickonsubmit" row.
[void]$wsEdit.Range("A7").Select()

# Switch the active tab from NewLoanInput to Edit Repayment Schedule.
$wsEdit.Activate()
